# Update cryptos list data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '38.344.13'
Set-TextValue "E2" '  +1.37%  '

Set-TextValue "D3" '2.102.38'
Set-TextValue "E3" '  +3.23%  '

Set-TextValue "E4" '  -0.33%  '

Set-TextValue "D5" '229.13'
Set-TextValue "E5" '  +0.43%  '

Set-TextValue "E6" '  +1.48%  '

Set-TextValue "D7" '61.25'
Set-TextValue "E7" '  +1.62%  '

Set-TextValue "E8" '  -0.02%  '

Set-TextValue "E9" '  +0.55%  '

Set-TextValue "D10" '0.0846'
Set-TextValue "E10" '  +3.27%  '

Set-TextValue "E11" '  +0.61%  '

Set-TextValue "D12" '2.416.17'
Set-TextValue "E12" '  +3.26%  '

Set-TextValue "D13" '14.76'
Set-TextValue "E13" '  +1.92%  '

Set-TextValue "D14" '22.36'
Set-TextValue "E14" '  +5.89%  '

Set-TextValue "E15" '  +6.02%  '

Set-TextValue "E16" '  +2.55%  '

Set-TextValue "D17" '2.103.08'
Set-TextValue "E17" '  +3.20%  '

Set-TextValue "D18" '38.269.23'
Set-TextValue "E18" '  +1.24%  '

Set-TextValue "D19" '6.03'
Set-TextValue "E19" '  +2.11%  '

Set-TextValue "D20" '70.42'
Set-TextValue "E20" '  +0.89%  '

Set-TextValue "E21" '  +1.39%  '

Set-TextValue "E22" '  +0.29%  '

Set-TextValue "E23" '  +0.02%  '

Set-TextValue "E24" '  +0.72%  '

Set-TextValue "D25" '2.32'
Set-TextValue "E25" '  +3.01%  '

Set-TextValue "D26" '169.82'
Set-TextValue "E26" '  +1.29%  '

Set-TextValue "E27" '  +0.47%  '

Set-TextValue "E28" '  +0.78%  '

Set-TextValue "D29" '19.05'
Set-TextValue "E29" '  +0.96%  '

Set-TextValue "E30" '  +6.25%  '

Set-TextValue "E31" '  -0.86%  '

Set-TextValue "D32" '2.38'
Set-TextValue "E32" '  +8.08%  '

Set-TextValue "E33" '  +4.75%  '

Set-TextValue "D34" '4.44'
Set-TextValue "E34" '  +0.87%  '

Set-TextValue "D35" '0.0606'
Set-TextValue "E35" '  -0.06%  '

Set-TextValue "B36" 'LidoDAOToken'
Set-TextValue "C36" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D36" '2.39'
Set-TextValue "E36" '  +4.95%  '

Set-TextValue "B37" 'THORChain'
Set-TextValue "C37" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D37" '6.42'
Set-TextValue "E37" '  +1.51%  '

Set-TextValue "E38" '  +5.67%  '

Set-TextValue "E39" '  +0.02%  '

Set-TextValue "D40" '18.11'
Set-TextValue "E40" '  +1.99%  '

Set-TextValue "D41" '1.548.72'
Set-TextValue "E41" '  +0.79%  '

Set-TextValue "D42" '100.08'
Set-TextValue "E42" '  +3.86%  '

Set-TextValue "D43" '0.0219'
Set-TextValue "E43" '  +0.96%  '

Set-TextValue "E44" '  +1.60%  '

Set-TextValue "D45" '0.0910'
Set-TextValue "E45" '  -0.33%  '

Set-TextValue "D46" '4.15'
Set-TextValue "E46" '  +3.51%  '

Set-TextValue "E47" '  +1.44%  '

Set-TextValue "D48" '7.50'
Set-TextValue "E48" '  +4.70%  '

Set-TextValue "E49" '  +1.95%  '

Set-TextValue "D51" '2.300.86'
Set-TextValue "E51" '  +3.21%  '
